# Insert a new row at position 768, which pushes the existing row 768
# (and everything below it) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(768).Insert()

# Populate the newly inserted row 768 with the new data record.
$ws.Range("A768").Value = 8
$ws.Range("B768").Value = "Terminal La Palmera de La Serena"
$ws.Range("C768").Value = "Coquimbo"
$ws.Range("D768").Value = 45124
$ws.Range("E768").Value = 4
$ws.Range("F768").Value = 100112024
$ws.Range("G768").Value = "Choclo"
$ws.Range("H768").Value = "Dulce o Americano"
$ws.Range("I768").Value = "Primera"
$ws.Range("J768").Value = 400
$ws.Range("K768").Value = 27000
$ws.Range("L768").Value = 28000
$ws.Range("M768").Value = 27500
$ws.Range("N768").Value = "$/malla 70 unidades"
$ws.Range("O768").Value = "Región de Arica y Parinacota"
$ws.Range("P768").Value = 393
$ws.Range("Q768").Value = 70
$ws.Range("R768").Value = "Hortaliza"
